$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$notes = $s.NotesPage
$shape = $notes.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$tr.InsertAfter([char]13 + "NOTEPARA1" + [char]13 + "NOTEPARA2")
